$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force these price cells to stay plain text (they are stored as
# text in the workbook even though they look numeric, e.g. "1.000"
# or "245.42" -- without this Excel would coerce them to numbers
# and normalize/round their display).
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '30.577.60'
$ws.Range('E2').Value = '  +0.01%  '
$ws.Range('D3').Value = '1.932.30'
$ws.Range('E3').Value = '  +0.49%  '
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '245.42'
$ws.Range('E5').Value = '  -0.84%  '
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  +0.01%  '
$ws.Range('D7').Value = '0.4843'
$ws.Range('E7').Value = '  +2.33%  '
$ws.Range('D8').Value = '0.2918'
$ws.Range('E8').Value = '  -0.14%  '
$ws.Range('D9').Value = '0.06790'
$ws.Range('E9').Value = '  -0.72%  '
$ws.Range('D10').Value = '113.00'
$ws.Range('E10').Value = '  +5.79%  '
$ws.Range('D11').Value = '19.42'
$ws.Range('E11').Value = '  +4.67%  '
$ws.Range('D12').Value = '1.939.68'
$ws.Range('E12').Value = '  +0.86%  '
$ws.Range('D13').Value = '0.07599'
$ws.Range('E13').Value = '  -1.59%  '
$ws.Range('D14').Value = '5.489'
$ws.Range('E14').Value = '  +2.83%  '
$ws.Range('D15').Value = '0.6800'
$ws.Range('E15').Value = '  +0.64%  '
$ws.Range('D16').Value = '298.59'
$ws.Range('E16').Value = '  +2.74%  '
$ws.Range('D17').Value = '30.591.28'
$ws.Range('E18').Value = '  +0.80%  '
$ws.Range('D19').Value = '0.000007653'
$ws.Range('E19').Value = '  +0.31%  '
$ws.Range('B20').Value = 'Dai'
$ws.Range('C20').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D20').Value = '0.9996'
$ws.Range('E20').Value = '  -0.06%  '
$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').Value = '2.188.46'
$ws.Range('E21').Value = '  +0.26%  '
$ws.Range('D22').Value = '5.545'
$ws.Range('E22').Value = '  -0.63%  '
$ws.Range('E23').Value = '  -0.12%  '
$ws.Range('D24').Value = '6.521'
$ws.Range('E24').Value = '  +0.39%  '
$ws.Range('D25').Value = '9.572'
$ws.Range('E25').Value = '  +0.17%  '
$ws.Range('D26').Value = '168.52'
$ws.Range('E26').Value = '  +0.92%  '
$ws.Range('D27').Value = '20.35'
$ws.Range('E27').Value = '  -2.19%  '
$ws.Range('D28').Value = '2.124'
$ws.Range('E28').Value = '  -0.37%  '
$ws.Range('D29').Value = '0.1070'
$ws.Range('E29').Value = '  -0.21%  '
$ws.Range('D30').Value = '1.430'
$ws.Range('E30').Value = '  +1.58%  '
$ws.Range('D31').Value = '4.190'
$ws.Range('E31').Value = '  -0.54%  '
$ws.Range('D32').Value = '4.100'
$ws.Range('E32').Value = '  +0.10%  '
$ws.Range('D33').Value = '0.05002'
$ws.Range('E33').Value = '  -1.24%  '
$ws.Range('D34').Value = '0.7506'
$ws.Range('E34').Value = '  +1.60%  '
$ws.Range('E35').Value = '  -0.13%  '
$ws.Range('D36').Value = '0.02040'
$ws.Range('E36').Value = '  -0.67%  '
$ws.Range('E37').Value = '  -0.91%  '
$ws.Range('D38').Value = '2.694'
$ws.Range('E38').Value = '  +0.49%  '
$ws.Range('D39').Value = '2.025'
$ws.Range('E39').Value = '  -1.59%  '
$ws.Range('D40').Value = '110.09'
$ws.Range('E40').Value = '  -1.32%  '
$ws.Range('D41').Value = '0.4470'
$ws.Range('E41').Value = '  -1.33%  '
$ws.Range('D42').Value = '0.8725'
$ws.Range('E42').Value = '  -0.43%  '
$ws.Range('E43').Value = '  -1.50%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').Value = '70.13'
$ws.Range('E44').Value = '  +3.00%  '
$ws.Range('B45').Value = 'PaxDollar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D45').Value = '1.001'
$ws.Range('E45').Value = '  +0.12%  '
$ws.Range('D46').Value = '7.337'
$ws.Range('E46').Value = '  -0.11%  '
$ws.Range('D47').Value = '49.34'
$ws.Range('E47').Value = '  +0.68%  '
$ws.Range('D48').Value = '9.355'
$ws.Range('E48').Value = '  -0.78%  '
$ws.Range('E49').Value = '  -2.97%  '
$ws.Range('D50').Value = '0.2546'
$ws.Range('E50').Value = '  +2.16%  '
$ws.Range('D51').Value = '35.11'
$ws.Range('E51').Value = '  -0.67%  '
